$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '323.71'
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '8.98%'
$c = $ws.Range("G2")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '49.59'
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '18.63%'
$c = $ws.Range("G3")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '5.357'
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '7.01%'
$c = $ws.Range("G4")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '0.08156'
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '8.40%'
$c = $ws.Range("G5")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '4.614'
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '5.42%'
$c = $ws.Range("G6")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '1.657'
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '4.68%'
$c = $ws.Range("G7")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '1.165'
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '25.69%'
$c = $ws.Range("G8")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.1347'
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '12.81%'
$c = $ws.Range("G9")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.1961'
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '7.52%'
$c = $ws.Range("G10")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.09547'
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '6.92%'
$c = $ws.Range("G11")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.04370'
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '7.07%'
$c = $ws.Range("G12")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.1048'
$c = $ws.Range("G13")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.001329'
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '3.44%'
$c = $ws.Range("G14")
$c.NumberFormat = "@"
$c.Value = '5'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.005959'
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '0.50%'
$c = $ws.Range("G15")
$c.NumberFormat = "@"
$c.Value = '5'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '3.393'
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '1.11%'
$c = $ws.Range("G16")
$c.NumberFormat = "@"
$c.Value = '5'
$ws.Range("B17").Value = 'BTSEToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '2.438'
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '1.52%'
$c = $ws.Range("G17")
$c.NumberFormat = "@"
$c.Value = '5'
$ws.Range("B18").Value = 'BitpandaEcosystemToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.3393'
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '2.39%'
$c = $ws.Range("G18")
$c.NumberFormat = "@"
$c.Value = '5'
$ws.Range("B19").Value = 'MCDex'
$ws.Range("C19").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '8.183'
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '0.50%'
$c = $ws.Range("G19")
$c.NumberFormat = "@"
$c.Value = '5'
$ws.Range("B20").Value = 'ProBitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.1420'
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '2.08%'
$c = $ws.Range("G20")
$c.NumberFormat = "@"
$c.Value = '5'
$ws.Range("B21").Value = 'ZBToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '0.3054'
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '-1.54%'
$c = $ws.Range("G21")
$c.NumberFormat = "@"
$c.Value = '5'
$ws.Range("B22").Value = 'CoinExToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.04316'
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '5.50%'
$c = $ws.Range("G22")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.001305'
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '3.06%'
$c = $ws.Range("G23")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '0.004264'
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '9.28%'
$c = $ws.Range("G24")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '9.72%'
$c = $ws.Range("G25")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '0.0003723'
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '-0.03%'
$c = $ws.Range("G26")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("G27")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("G28")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("G29")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("G30")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("G31")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("G32")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("G33")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("G34")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("G35")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("G36")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("G37")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.02788'
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '16.07%'
$c = $ws.Range("G38")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.05545'
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '5.99%'
$c = $ws.Range("G39")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.006201'
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '-1.65%'
$c = $ws.Range("G40")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.007739'
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '-1.01%'
$c = $ws.Range("G41")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.1451'
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '9.41%'
$c = $ws.Range("G42")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.007681'
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '3.74%'
$c = $ws.Range("G43")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '11.67%'
$c = $ws.Range("G44")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.3519'
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '18.69%'
$c = $ws.Range("G45")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.00006773'
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '2.83%'
$c = $ws.Range("G46")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '-0.06%'
$c = $ws.Range("G47")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '94.07%'
$c = $ws.Range("G48")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.004001'
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '-4.81%'
$c = $ws.Range("G49")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '-0.06%'
$c = $ws.Range("G50")
$c.NumberFormat = "@"
$c.Value = '5'
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '-0.06%'
$c = $ws.Range("G51")
$c.NumberFormat = "@"
$c.Value = '5'
